# Fruta / hortaliza, semanal
# A new weekly observation (date serial 44785, price 160) is inserted at the
# top of the date/price series that lives in rows 282-350 (columns D and J).
# This pushes every existing D/J pair down by one row; the pair that falls
# off the bottom (the old row 350 values) becomes a brand-new row 351, which
# is otherwise an exact copy of row 350.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Create the new row 351 as a full copy of row 350 (this captures row
#    350's original D/J values before they get overwritten by the shift
#    below).
$src = $ws.Range("A350:R350")
$dst = $ws.Range("A351:R351")
$src.Copy($dst)

# 2) Shift the D (date) and J (price) values down by one row, working from
#    the bottom (350) up to 283 so that each row's original value is read
#    before it gets overwritten.
for ($r = 350; $r -ge 283; $r--) {
    $prev = $r - 1
    $prevDate = $ws.Cells.Item($prev, 4).Value2
    $prevPrice = $ws.Cells.Item($prev, 10).Value2
    $ws.Cells.Item($r, 4).Value2 = $prevDate
    $ws.Cells.Item($r, 10).Value2 = $prevPrice
}

# 3) Insert the brand-new observation in row 282.
$ws.Cells.Item(282, 4).Value2 = 44785
$ws.Cells.Item(282, 10).Value2 = 160
